$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple per-row price / volume updates ---
$ws.Range("D2").Value = '67.711.28'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '2.524.17'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue $ws.Range("D5") '592.02'
$ws.Range("E5").Value = '  +1.96%  '
Set-TextValue $ws.Range("D6") '175.63'
$ws.Range("E6").Value = '  +4.63%  '
$ws.Range("E7").Value = '  +0.12%  '
Set-TextValue $ws.Range("D8") '0.529'
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").Value = '2.524.93'
$ws.Range("E9").Value = '  -0.46%  '
Set-TextValue $ws.Range("D10") '0.141'
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("E11").Value = '  +2.29%  '
Set-TextValue $ws.Range("D12") '5.15'
$ws.Range("E12").Value = '  +0.35%  '
Set-TextValue $ws.Range("D13") '0.344'
$ws.Range("E13").Value = '  -1.86%  '
Set-TextValue $ws.Range("D14") '26.78'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '2.983.46'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '67.541.44'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '2.526.56'
$ws.Range("E18").Value = '  -0.53%  '
Set-TextValue $ws.Range("D19") '8.01'
$ws.Range("E19").Value = '  +5.01%  '
Set-TextValue $ws.Range("D20") '11.42'
$ws.Range("E20").Value = '  +0.92%  '
Set-TextValue $ws.Range("D21") '360.33'
$ws.Range("E21").Value = '  +3.91%  '
Set-TextValue $ws.Range("D22") '4.20'
$ws.Range("E22").Value = '  +0.29%  '
Set-TextValue $ws.Range("D23") '4.64'
$ws.Range("E23").Value = '  +1.91%  '
Set-TextValue $ws.Range("D24") '1.97'
$ws.Range("E24").Value = '  +2.14%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D30").Value = '0.0₃0986'
$ws.Range("E30").Value = '  +0.87%  '
Set-TextValue $ws.Range("D31") '551.98'
$ws.Range("E31").Value = '  +5.15%  '
Set-TextValue $ws.Range("D32") '8.24'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("E33").Value = '  +2.57%  '
Set-TextValue $ws.Range("D34") '1.86'
$ws.Range("E34").Value = '  +2.65%  '
$ws.Range("E35").Value = '  +0.13%  '
Set-TextValue $ws.Range("D36") '1.00'
$ws.Range("E36").Value = '  +0.11%  '
Set-TextValue $ws.Range("D37") '1.47'
$ws.Range("E37").Value = '  +1.61%  '
Set-TextValue $ws.Range("D38") '155.87'
$ws.Range("E38").Value = '  -0.61%  '
Set-TextValue $ws.Range("D39") '18.74'
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("E40").Value = '  +1.83%  '
Set-TextValue $ws.Range("D41") '0.356'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  +2.88%  '
Set-TextValue $ws.Range("D43") '5.17'
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("D48").Value = '0.0₆0280'
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("E50").Value = '  -0.95%  '
Set-TextValue $ws.Range("D51") '0.0756'
$ws.Range("E51").Value = '  +0.06%  '

# --- Rows whose Coin/Link/Price/Volume data were reordered (rank swap) ---
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D26") '70.88'
$ws.Range("E26").Value = '  +2.77%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D27") '10.28'
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D44") '1.00'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D45") '2.51'
$ws.Range("E45").Value = '  +4.10%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D46") '147.07'
$ws.Range("E46").Value = '  -0.27%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D47") '0.561'
$ws.Range("E47").Value = '  +0.85%  '
